{"js": "// Update the two-digit x two-digit multiplication \"answers\" table.\n// Each data cell holds a single run of plain text like \"49\u00d766=3234\".\n// We replace the text of each matching cell in place (Replace insert\n// location) so the existing run/paragraph formatting (font, size,\n// alignment) is preserved exactly, and we only rewrite the <w:t> value.\n\nconst replacements = [\n  [\"49\u00d766=3234\", \"20\u00d741=820\"],\n  [\"36\u00d763=2268\", \"46\u00d780=3680\"],\n  [\"33\u00d745=1485\", \"25\u00d771=1775\"],\n  [\"23\u00d770=1610\", \"48\u00d730=1440\"],\n  [\"89\u00d767=5963\", \"58\u00d799=5742\"],\n  [\"84\u00d766=5544\", \"71\u00d732=2272\"],\n  [\"31\u00d718=558\", \"30\u00d724=720\"],\n  [\"25\u00d726=650\", \"94\u00d728=2632\"],\n  [\"31\u00d762=1922\", \"94\u00d770=6580\"],\n  [\"50\u00d722=1100\", \"95\u00d761=5795\"],\n  [\"82\u00d796=7872\", \"79\u00d741=3239\"],\n  [\"94\u00d788=8272\", \"45\u00d783=3735\"],\n  [\"40\u00d724=960\", \"83\u00d722=1826\"],\n  [\"97\u00d747=4559\", \"32\u00d781=2592\"],\n  [\"16\u00d735=560\", \"75\u00d721=1575\"],\n  [\"51\u00d730=1530\", \"52\u00d782=4264\"],\n  [\"46\u00d723=1058\", \"94\u00d763=5922\"],\n  [\"67\u00d786=5762\", \"81\u00d789=7209\"],\n  [\"44\u00d780=3520\", \"16\u00d735=560\"],\n  [\"49\u00d769=3381\", \"75\u00d754=4050\"],\n  [\"74\u00d758=4292\", \"28\u00d755=1540\"],\n  [\"86\u00d744=3784\", \"13\u00d788=1144\"],\n  [\"55\u00d757=3135\", \"16\u00d780=1280\"],\n  [\"97\u00d779=7663\", \"79\u00d791=7189\"],\n  [\"57\u00d751=2907\", \"61\u00d737=2257\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Gather every table cell's body (there is a single table in this\n// document, but walk every table defensively in case that changes).\nconst cellBodies = [];\nfor (let t = 0; t < tables.items.length; t++) {\n  const table = tables.items[t];\n  table.rows.load(\"items\");\n}\nawait context.sync();\n\nfor (let t = 0; t < tables.items.length; t++) {\n  const rows = tables.items[t].rows.items;\n  for (let r = 0; r < rows.length; r++) {\n    rows[r].cells.load(\"items\");\n  }\n}\nawait context.sync();\n\nfor (let t = 0; t < tables.items.length; t++) {\n  const rows = tables.items[t].rows.items;\n  for (let r = 0; r < rows.length; r++) {\n    const cells = rows[r].cells.items;\n    for (let c = 0; c < cells.length; c++) {\n      cellBodies.push(cells[c].body);\n    }\n  }\n}\n\n// Snapshot each cell's current text exactly once, BEFORE any edits, so\n// that a new value which happens to equal another cell's old value\n// (e.g. \"16\u00d735=560\" is both an old value in one cell and the new value\n// written into a different cell) can never be double-matched.\nfor (const body of cellBodies) {\n  body.paragraphs.load(\"items\");\n}\nawait context.sync();\n\nconst targets = [];\nfor (const body of cellBodies) {\n  const paras = body.paragraphs.items;\n  if (paras.length === 0) continue;\n  const para = paras[0];\n  para.load(\"text\");\n  targets.push(para);\n}\nawait context.sync();\n\nfor (const para of targets) {\n  const current = para.text;\n  for (const [oldText, newText] of replacements) {\n    if (current === oldText) {\n      para.insertText(newText, Word.InsertLocation.replace);\n      break;\n    }\n  }\n}\nawait context.sync();\n", "ps1": "# Update the two-digit x two-digit multiplication \"answers\" table.\n# Each data cell holds a single run of plain text like \"49\u00d766=3234\".\n# We overwrite Cell.Range.Text in place, which preserves the existing\n# run/paragraph formatting (font, size, alignment) and only rewrites\n# the text content, matching the canonical diff exactly.\n\n$d = $word.ActiveDocument\n\n$oldValues = @(\n  \"49\u00d766=3234\", \"36\u00d763=2268\", \"33\u00d745=1485\", \"23\u00d770=1610\", \"89\u00d767=5963\",\n  \"84\u00d766=5544\", \"31\u00d718=558\",  \"25\u00d726=650\",  \"31\u00d762=1922\", \"50\u00d722=1100\",\n  \"82\u00d796=7872\", \"94\u00d788=8272\", \"40\u00d724=960\",  \"97\u00d747=4559\", \"16\u00d735=560\",\n  \"51\u00d730=1530\", \"46\u00d723=1058\", \"67\u00d786=5762\", \"44\u00d780=3520\", \"49\u00d769=3381\",\n  \"74\u00d758=4292\", \"86\u00d744=3784\", \"55\u00d757=3135\", \"97\u00d779=7663\", \"57\u00d751=2907\"\n)\n\n$newValues = @(\n  \"20\u00d741=820\",  \"46\u00d780=3680\", \"25\u00d771=1775\", \"48\u00d730=1440\", \"58\u00d799=5742\",\n  \"71\u00d732=2272\", \"30\u00d724=720\",  \"94\u00d728=2632\", \"94\u00d770=6580\", \"95\u00d761=5795\",\n  \"79\u00d741=3239\", \"45\u00d783=3735\", \"83\u00d722=1826\", \"32\u00d781=2592\", \"75\u00d721=1575\",\n  \"52\u00d782=4264\", \"94\u00d763=5922\", \"81\u00d789=7209\", \"16\u00d735=560\",  \"75\u00d754=4050\",\n  \"28\u00d755=1540\", \"13\u00d788=1144\", \"16\u00d780=1280\", \"79\u00d791=7189\", \"61\u00d737=2257\"\n)\n\n# Snapshot every table cell's current text BEFORE making any edits. This\n# guarantees the lookup below always matches against the document's\n# original content, even though some of the new values happen to equal\n# another cell's original value (e.g. \"16\u00d735=560\" is an old value in one\n# cell and the new value written into a different cell). Walk every\n# table defensively (this document has a single table, with data rows\n# separated by blank rows, but a cell-text match does not depend on\n# that layout).\n$cells = New-Object System.Collections.ArrayList\n$originalTexts = New-Object System.Collections.ArrayList\nforeach ($tbl in $d.Tables) {\n  for ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {\n      $cell = $null\n      try {\n        $cell = $tbl.Cell($r, $c)\n      } catch {\n        $cell = $null\n      }\n      if ($cell -ne $null) {\n        $text = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        [void]$cells.Add($cell)\n        [void]$originalTexts.Add($text)\n      }\n    }\n  }\n}\n\nfor ($i = 0; $i -lt $cells.Count; $i++) {\n  $text = $originalTexts[$i]\n  for ($j = 0; $j -lt $oldValues.Length; $j++) {\n    if ($text -eq $oldValues[$j]) {\n      $cells[$i].Range.Text = $newValues[$j]\n      break\n    }\n  }\n}\n"}
